$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New "Performance Targets" block (rows 25-32), literal values first so the
#    formulas calculate against real numbers.
# ---------------------------------------------------------------------------

# Row 25 header text
$ws.Range("A25").Value = "Performance Targets (Top Students from 2016)"

# Row labels (column A)
$ws.Range("A26").Value = "TARGET #1 (Delay) [Samuel Bauza & Michael Petrucci] "
$ws.Range("A27").Value = "TARGET #2 (Delay) [Jinzhe Zhang]"
$ws.Range("A28").Value = "TARGET #3 (Delay) [Kevin Pan & Brandon Ustaris]"
$ws.Range("A30").Value = "TARGET #1 (Area x Delay) [Jinzhe Zhang]"
$ws.Range("A31").Value = "TARGET #2 (Area x Delay) [Kevin Pan & Brandon Ustaris]"
$ws.Range("A32").Value = "TARGET #3 (Area x Delay) [Rui Han]"

# Literal numeric inputs: B (ALUTS), C (Registers), E (FMAX), G (#Cycles)
$ws.Range("B26").Value = 1319
$ws.Range("C26").Value = 1508
$ws.Range("E26").Value = 186.15
$ws.Range("G26").Value = 146

$ws.Range("B27").Value = 594
$ws.Range("C27").Value = 968
$ws.Range("E27").Value = 321.23
$ws.Range("G27").Value = 252

$ws.Range("B28").Value = 559
$ws.Range("C28").Value = 1007
$ws.Range("E28").Value = 325.1
$ws.Range("G28").Value = 258

$ws.Range("B30").Value = 594
$ws.Range("C30").Value = 968
$ws.Range("E30").Value = 321.23
$ws.Range("G30").Value = 252

$ws.Range("B31").Value = 559
$ws.Range("C31").Value = 1007
$ws.Range("E31").Value = 325.1
$ws.Range("G31").Value = 258

$ws.Range("B32").Value = 786
$ws.Range("C32").Value = 897
$ws.Range("E32").Value = 293.17
$ws.Range("G32").Value = 249

# Formulas (column D, F, H, I) - set across the whole block at once so Excel
# builds one shared-formula group per column, same as a fill-down would.
$ws.Range("D26:D32").Formula = "=B26+C26"
$ws.Range("F26:F32").Formula = "=(1/E26)*1000"
$ws.Range("H26:H32").Formula = "=F26*G26"
$ws.Range("I26:I32").Formula = "=D26*H26*(10^-9)"

# ---------------------------------------------------------------------------
# 2. Formatting - build up fills/borders/number formats in the same order the
#    original author's workbook encodes them (xf indices 4-9).
# ---------------------------------------------------------------------------

# xf4: yellow fill, no border -> core data cells B:H on rows 26-32
$ws.Range("B26:H32").Interior.Color = 65535

# xf5: yellow fill + scientific number format -> I26:I32
$ws.Range("I26:I32").Interior.Color = 65535
$ws.Range("I26:I32").NumberFormat = "0.00E+00"

# xf6: yellow fill + bottom border -> header row, B25:H25
$ws.Range("B25:H25").Interior.Color = 65535
$ws.Range("B25:H25").Borders(9).LineStyle = 1

# xf7: yellow fill + bottom border + scientific number format -> I25
$ws.Range("I25").Interior.Color = 65535
$ws.Range("I25").NumberFormat = "0.00E+00"
$ws.Range("I25").Borders(9).LineStyle = 1

# xf8: yellow fill + right & bottom border -> A25 (right first, then bottom,
# mirrors how the source file grew its border list: right-only before
# right+bottom)
$ws.Range("A25").Interior.Color = 65535
$ws.Range("A25").Borders(10).LineStyle = 1
$ws.Range("A25").Borders(9).LineStyle = 1

# xf9: yellow fill + right border only -> A26:A32
$ws.Range("A26:A32").Interior.Color = 65535
$ws.Range("A26:A32").Borders(10).LineStyle = 1

# Scientific-notation placeholder formatting that already existed on I2 now
# extends all the way down the sheet (rows 3-39, skipping the yellow block
# which already carries its own format).
$ws.Range("I3:I24").NumberFormat = "0.00E+00"
$ws.Range("I33:I39").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# 3. Column width + "Increment 0" label tweak (A2). Re-labelled last so the
#    shared-string table keeps the original low indices for the still-used
#    strings and simply appends the new label at the end.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 49.5703125
$ws.Range("A2").Value = "Increment 0 (starting point) [Winter 2016]"

# ---------------------------------------------------------------------------
# 4. Selection, matching the saved cursor position in the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("A5").Select()
